$d = $word.ActiveDocument

# Locate the target list-paragraph that currently reads "Vault operator "
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd([char]13, [char]7) -eq "Vault operator ") {
        $target = $cand
        $targetIndex = $i
        break
    }
}

$url = "https://stackoverflow.com/questions/42944115/is-there-a-way-to-run-hashicorp-vault-as-a-windows-service"
$placeholder = "@@LINKPLACEHOLDER@@"

$r = $target.Range
# Replace the run text, add a manual line break ([char]11 = vertical-tab = w:br),
# and leave a placeholder for the hyperlink's display text.
$r.Text = "Running the vault as service" + [char]11 + $placeholder

# Recompute the paragraph range after the text change, then select just the
# placeholder text (it sits right before the paragraph mark).
$full = $target.Range
$placeholderEnd = $full.End - 1
$placeholderStart = $placeholderEnd - $placeholder.Length
$linkRange = $d.Range($placeholderStart, $placeholderEnd)
$d.Hyperlinks.Add($linkRange, $url, $null, $null, $url) | Out-Null

# Insert a new, empty "List Paragraph"-styled paragraph right after the
# paragraph we just edited.
$following = $d.Paragraphs.Item($targetIndex + 1)
$following.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Style = "List Paragraph"
